$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 97.53587864304308
$ws.Cells.Item(3, 3).Value = 97.72309972103312
$ws.Cells.Item(4, 3).Value = 98.225715909605
$ws.Cells.Item(5, 3).Value = 99.23215930283021
$ws.Cells.Item(6, 3).Value = 99.60574330297312
$ws.Cells.Item(7, 3).Value = 99.7311943377388
$ws.Cells.Item(8, 3).Value = 99.23885965288227
$ws.Cells.Item(9, 3).Value = 99.19316953925974
$ws.Cells.Item(10, 3).Value = 98.96896067432068
$ws.Cells.Item(11, 3).Value = 99.50487987522162
$ws.Cells.Item(12, 3).Value = 99.85810390725479
$ws.Cells.Item(13, 3).Value = 100.0197159098832
$ws.Cells.Item(14, 3).Value = 99.24823068958128
$ws.Cells.Item(15, 3).Value = 99.38240424270809
$ws.Cells.Item(16, 3).Value = 99.30996303454629
$ws.Cells.Item(17, 3).Value = 99.77149822907566
$ws.Cells.Item(18, 3).Value = 99.99389333428415
$ws.Cells.Item(19, 3).Value = 100.5782742046739
$ws.Cells.Item(20, 3).Value = 99.87374509615711
$ws.Cells.Item(21, 3).Value = 99.90784354954091
$ws.Cells.Item(22, 3).Value = 99.79432607285075
$ws.Cells.Item(23, 3).Value = 100.2007512771496
$ws.Cells.Item(24, 3).Value = 100.3275046508865
$ws.Cells.Item(25, 3).Value = 100.2636237852996
$ws.Cells.Item(27, 3).Value = 99.29479070918457
$ws.Cells.Item(28, 3).Value = 100.0058691952096
$ws.Cells.Item(29, 3).Value = 100.8800542579492
$ws.Cells.Item(30, 3).Value = 101.2321686900629
$ws.Cells.Item(31, 3).Value = 101.5890063511758
$ws.Cells.Item(32, 3).Value = 100.7485646269179
$ws.Cells.Item(33, 3).Value = 101.0897405626647
$ws.Cells.Item(34, 3).Value = 101.4493300510087
$ws.Cells.Item(35, 3).Value = 102.8488341535093
$ws.Cells.Item(36, 3).Value = 103.0847344860595
$ws.Cells.Item(37, 3).Value = 104.0109908196817
$ws.Cells.Item(38, 3).Value = 103.9928503446066
$ws.Cells.Item(39, 3).Value = 104.6119230335139
$ws.Cells.Item(40, 3).Value = 107.2463749641436
$ws.Cells.Item(41, 3).Value = 107.09908649001
$ws.Cells.Item(42, 3).Value = 108.03923983085
$ws.Cells.Item(43, 3).Value = 110.0951765946902
$ws.Cells.Item(44, 3).Value = 110.1650409567812
$ws.Cells.Item(45, 3).Value = 110.3702972896204
$ws.Cells.Item(46, 3).Value = 109.5129602890809
$ws.Cells.Item(47, 3).Value = 109.7829129283668
$ws.Cells.Item(48, 3).Value = 109.9631740797092
$ws.Cells.Item(49, 3).Value = 110.3642651472358
$ws.Cells.Item(50, 3).Value = 110.7008284302895
$ws.Cells.Item(51, 3).Value = 111.4923326104787
$ws.Cells.Item(52, 3).Value = 112.0597714648144
$ws.Cells.Item(53, 3).Value = 113.0855831092472
$ws.Cells.Item(54, 3).Value = 113.2985267206117
$ws.Cells.Item(55, 3).Value = 114.3056822378065
$ws.Cells.Item(56, 3).Value = 114.707732463765
$ws.Cells.Item(57, 3).Value = 115.1718372274985
$ws.Cells.Item(58, 3).Value = 115.0606003567965
$ws.Cells.Item(59, 3).Value = 115.3089675665149
$ws.Cells.Item(60, 3).Value = 115.2586452857236
$ws.Cells.Item(61, 3).Value = 115.798008647439
$ws.Cells.Item(62, 3).Value = 116.6268333990807
$ws.Cells.Item(63, 3).Value = 117.4395330853832
$ws.Cells.Item(64, 3).Value = 118.3596043566196
$ws.Cells.Item(65, 3).Value = 118.9728503232974
$ws.Cells.Item(66, 3).Value = 120.4684466591
$ws.Cells.Item(67, 3).Value = 121.3091606668395
$ws.Cells.Item(68, 3).Value = 120.1155417966271
$ws.Cells.Item(69, 3).Value = 120.1480577386248
$ws.Cells.Item(70, 3).Value = 119.3461789259675
$ws.Cells.Item(71, 3).Value = 119.8976934087258
$ws.Cells.Item(72, 3).Value = 120.8714745717841
$ws.Cells.Item(73, 3).Value = 121.7818813028505
$ws.Cells.Item(74, 3).Value = 122.5943252861004
$ws.Cells.Item(75, 3).Value = 123.3985225672302
$ws.Cells.Item(76, 3).Value = 123.4463540894337
$ws.Cells.Item(77, 3).Value = 124.6970524631508
$ws.Cells.Item(78, 3).Value = 125.1928753531599
$ws.Cells.Item(79, 3).Value = 126.0525321005129
$ws.Cells.Item(80, 3).Value = 125.973386216252
